$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A52").Value = 44365.76761002894

$ws.Range("A53").Value = 44366.7671924891
$ws.Range("A53").NumberFormat = $ws.Range("A52").NumberFormat

$ws.Range("B53").Value = 78490
$ws.Range("C53").Value = 65989
$ws.Range("D53").Value = 3442
$ws.Range("E53").Value = 2112
$ws.Range("F53").Value = 1493
$ws.Range("G53").Value = 20778
$ws.Range("H53").Value = 1438
$ws.Range("I53").Value = 899
$ws.Range("J53").Value = 182
